# Add a new "2022-Q3" sheet (fund holdings) right after "总计" and before
# "2022-Q2", renumbering nothing else (Excel shifts sheetIds/positions for
# us), then insert the corresponding summary row into "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet before the current "2022-Q2" sheet
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)          # currently "2022-Q2"
$newSheet = $wb.Worksheets.Add($anchor)
$newSheet.Name = "2022-Q3"

# Reuse look & feel (sheetPr / pageMargins / header style) from the sheet
# that used to be "2022-Q2" (now pushed one slot to the right) so the new
# tab matches its siblings.
$template = $wb.Worksheets.Item(3)        # "2022-Q2" (template layout)
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$template.Range("A2:H2").Copy()
$newSheet.Range("A2:H7").PasteSpecial(-4122)   # xlPasteFormats

$newSheet.Outline.SummaryBelow = $true
$newSheet.Outline.SummaryRight = $true
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings data (code, name, size, stock-position, position-ratio,
# market value, position rank). Columns B-G are stored as text in the
# source workbook (to avoid losing trailing zeros / leading zeros), H is
# numeric.
$rows = @(
    @("010114", "华宝新兴成长混合", "3.18", "80.50", "3.66", "0.1164", 10),
    @("310368", "申万菱信竞争优势混合A", "1.05", "92.76", "4.63", "0.0486", 7),
    @("004320", "前海开源沪港深乐享生活灵活配置混合", "0.27", "71.16", "3.29", "0.0089", 7),
    @("007463", "东海科技动力混合C", "0.13", "81.12", "5.44", "0.0071", 5),
    @("007439", "东海科技动力混合A", "0.12", "81.12", "5.44", "0.0065", 5),
    @("015173", "申万菱信竞争优势混合C", "0.11", "92.76", "4.63", "0.0051", 7)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2            # A: 0-based running index
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]     # B: 基金代码 (text, keep leading zeros)
    $newSheet.Cells.Item($r, 3).Value = $row[1]           # C: 基金名称
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]     # D: 基金规模 (text)
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]     # E: 股票总仓位 (text)
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]     # F: 仓位占比 (text)
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]     # G: 持有市值(亿元) (text)
    $newSheet.Cells.Item($r, 8).Value = $row[6]           # H: 仓位排名 (number)
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new leading data row for 2022-Q3
#    and shift the rest down (values rewritten directly, row by row).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)         # "总计"

# Extend the index-column formatting (style of A2:A7) down to the new A8.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

$summaryRows = @(
    @("2022-Q3", 6, 0.19),
    @("2022-Q2", 19, 2.75),
    @("2022-Q1", 9, 1.98),
    @("2021-Q4", 7, 5.06),
    @("2021-Q2", 4, 1.11),
    @("2021-Q1", 2, 0.02),
    @("2020-Q4", 2, 2.57)
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $r - 2     # A: running index
    $summary.Cells.Item($r, 2).Value = $row[0]    # B: 日期
    $summary.Cells.Item($r, 3).Value = $row[1]    # C: 持有数量(只)
    $summary.Cells.Item($r, 4).Value = $row[2]    # D: 持有市值(亿元)
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Keep the original active tab ("2020-Q4", the last sheet) selected,
#    since adding a worksheet moves the selection to the new tab.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()

